$p = $ppt.ActivePresentation

# Add the new slide by duplicating the existing "Title and Content" slide
# (slide 2) and dropping it right after it (slide 3). Duplicating (rather
# than Slides.Add) carries over the boilerplate PowerPoint always writes for
# an authored slide (clrMapOvr, group xfrm, slide-level creationId, etc.)
# instead of the bare-bones shell a freshly inserted layout gets.
$src = $p.Slides.Item(2)
$s = $src.Duplicate().Item(1)

# Drop the two picture shapes that came along with the duplicated slide -
# this new slide is text-only.
$s.Shapes.Item(4).Delete()
$s.Shapes.Item(3).Delete()

$s.Shapes.Item(1).Name = "Title 2"
$s.Shapes.Item(2).Name = "Content Placeholder 3"

# --- Title placeholder -----------------------------------------------------
$title = $s.Shapes.Item(1).TextFrame.TextRange
$title.Text = "Key Concepts of PyTorch"

# --- Content placeholder ---------------------------------------------------
# Build the bulleted definition list paragraph-by-paragraph (rather than one
# multi-line Text assignment) and only insert the plain-text runs first,
# applying Bold afterwards to the lead term of each line. Doing it in this
# order keeps every run's rPr/lang attribute intact.
$body = $s.Shapes.Item(2).TextFrame.TextRange

$body.Text = "Tensor"
[void]$body.InsertAfter(": Multi-dimensional array of numbers")

[void]$body.InsertAfter("`rAutograd")
[void]$body.InsertAfter(": A library for automatically computing gradients")

[void]$body.InsertAfter("`rNeural network module")
[void]$body.InsertAfter(": Base class for all neural network classes")

[void]$body.InsertAfter("`rOptimizer")
[void]$body.InsertAfter(": An algorithm that is used to adjust the parameters of a neural network in order to minimize a loss function")

$body.Paragraphs(1).Characters(1, 6).Font.Bold = $true
$body.Paragraphs(2).Characters(1, 8).Font.Bold = $true
$body.Paragraphs(3).Characters(1, 21).Font.Bold = $true
$body.Paragraphs(4).Characters(1, 9).Font.Bold = $true
